# Sheet1 holds a monthly time series (date in column A, two price-index
# values in columns B/C). This update:
#   1. Reorders the existing 2018-2021 data so each year's Oct/Nov/Dec rows
#      come first, followed by Jan..Sep of that year.
#   2. Appends new 2022 (Oct/Nov/Dec, then Jan..Sep) and 2023 (Jan..Jul) data.
# The full, correctly ordered A:B:C data set (67 rows covering 2018-10
# through 2023-07) is written out below, rewriting A2:C68 in one pass.
$data = @(
    ,@("2018-10", 100.5, 101.7)
    ,@("2018-11", 100.9, 101.6)
    ,@("2018-12", 101.4, 101.5)
    ,@("2018-01", 99.71899999999999, 102.2401)
    ,@("2018-02", 99.3, 102.1)
    ,@("2018-03", 99.3, 102.2)
    ,@("2018-04", 99.59999999999999, 101)
    ,@("2018-05", 100, 100.9)
    ,@("2018-06", 100.2, 101.1)
    ,@("2018-07", 100.4, 101.2)
    ,@("2018-08", 100.2, 101.5)
    ,@("2018-09", 100.5, 102)
    ,@("2019-10", 100.8, 100.8)
    ,@("2019-11", 100.9, 100.5)
    ,@("2019-12", 100.3, 100.4)
    ,@("2019-01", 102.7, 101.5)
    ,@("2019-02", 102.5, 101.6)
    ,@("2019-03", 102.4, 100.7)
    ,@("2019-04", 102, 101.6)
    ,@("2019-05", 102.3, 101.6)
    ,@("2019-06", 102.7, 101.6)
    ,@("2019-07", 102.1, 101.4)
    ,@("2019-08", 101.6, 101)
    ,@("2019-09", 101.8, 101)
    ,@("2020-10", 100.2, 100.4)
    ,@("2020-11", 100, 100.8)
    ,@("2020-12", 99.90000000000001, 100.9)
    ,@("2020-01", 98.7, 100.3)
    ,@("2020-02", 99, 100.1)
    ,@("2020-03", 99.09999999999999, 100.8)
    ,@("2020-04", 99.3, 100.7)
    ,@("2020-05", 98.59999999999999, 100.7)
    ,@("2020-06", 98.2, 100.9)
    ,@("2020-07", 98.5, 100.9)
    ,@("2020-08", 99.90000000000001, 100.9)
    ,@("2020-09", 99.59999999999999, 100.3)
    ,@("2021-10", 100, 101.2)
    ,@("2021-11", 102.3, 101.2)
    ,@("2021-12", 102.2, 101.2)
    ,@("2021-01", 101.9, 100.4)
    ,@("2021-02", 100.3, 99.90000000000001)
    ,@("2021-03", 100.2, 99.7)
    ,@("2021-04", 101.1, 99.90000000000001)
    ,@("2021-05", 102.7, 99.7)
    ,@("2021-06", 100.5, 100)
    ,@("2021-07", 101, 100.4)
    ,@("2021-08", 101.1, 100.5)
    ,@("2021-09", 100.5, 101)
    ,@("2022-10", 104.1, 102.4)
    ,@("2022-11", 102, 102.1)
    ,@("2022-12", 101.8, 101.9)
    ,@("2022-01", 101.6, 101.5)
    ,@("2022-02", 103.1, 101.6)
    ,@("2022-03", 103.2, 102.1)
    ,@("2022-04", 102.1, 101.9)
    ,@("2022-05", 101.1, 102.3)
    ,@("2022-06", 102.8, 102.6)
    ,@("2022-07", 102.8, 102)
    ,@("2022-08", 102.4, 102)
    ,@("2022-09", 103.1, 102.1)
    ,@("2023-01", 100.7, 101.5)
    ,@("2023-02", 99.09999999999999, 101.3)
    ,@("2023-03", 98.90000000000001, 100.8)
    ,@("2023-04", 98.59999999999999, 100.6)
    ,@("2023-05", 98.09999999999999, 100.2)
    ,@("2023-06", 98.09999999999999, 99.59999999999999)
    ,@("2023-07", 98.2, 99.59999999999999)
)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rowCount = $data.Count
$lastRow = 1 + $rowCount
$origLastRow = 49

for ($i = 0; $i -lt $rowCount; $i++) {
    $r = $i + 2
    $row = $data[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
}

# Newly added rows (beyond the original 49 data rows) need the same
# formatting (centered, bold, bordered) that column A already carries
# on the pre-existing rows. Copy that formatting down onto the new rows.
if ($lastRow -gt $origLastRow) {
    $ws.Range("A2").Copy()
    $ws.Range("A" + ($origLastRow + 1) + ":A" + $lastRow).PasteSpecial(-4122)
    $excel.CutCopyMode = 0
}

